# Act greficos y tablas web pob
# Rename sheets, reverse the Data series (Afrodescendiente / No afrodescendiente)
# and rewrite the Metadata ("Ficha tecnica") sheet with lowercase field-name
# keys, an extra "observaciones" row and a final "Mirador DESCA..." credit row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename worksheets: "Datos" -> "Data", "Ficha técnica" -> "Metadata"
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsMeta = $wb.Worksheets.Item(2)
$wsData.Name = "Data"
$wsMeta.Name = "Metadata"

# ---------------------------------------------------------------------------
# 2) Data sheet: rows 2-10 (years 2011..2019) are flipped top-to-bottom, so
#    the table now runs 2019 down to 2011 (Fecha, Afrodescendiente and
#    No afrodescendiente all move together as whole rows).
# ---------------------------------------------------------------------------
$wsData.Range("A2").Value = "2019"
$wsData.Range("B2").Value = 74.6
$wsData.Range("C2").Value = 62.9

$wsData.Range("A3").Value = "2018"
$wsData.Range("B3").Value = 76.4
$wsData.Range("C3").Value = 64

$wsData.Range("A4").Value = "2017"
$wsData.Range("B4").Value = 77.6
$wsData.Range("C4").Value = 64.2

$wsData.Range("A5").Value = "2016"
$wsData.Range("B5").Value = 79.6
$wsData.Range("C5").Value = 65.3

$wsData.Range("A6").Value = "2015"
$wsData.Range("B6").Value = 80.4
$wsData.Range("C6").Value = 65.5

$wsData.Range("A7").Value = "2014"
$wsData.Range("B7").Value = 80.4
$wsData.Range("C7").Value = 65.4

$wsData.Range("A8").Value = "2013"
$wsData.Range("B8").Value = 81
$wsData.Range("C8").Value = 66.7

$wsData.Range("A9").Value = "2012"
$wsData.Range("B9").Value = 81.5
$wsData.Range("C9").Value = 66.6

$wsData.Range("A10").Value = "2011"
$wsData.Range("B10").Value = 83
$wsData.Range("C10").Value = 67.5

# ---------------------------------------------------------------------------
# 3) Metadata sheet: keys switch from UPPERCASE labels to lowercase field
#    names, rows are reshuffled, a new "observaciones" row is inserted and
#    a "cita"/credit pair is appended at the end.
#    Row 1 (A1 blank / B1 " ") is unchanged, so it is left untouched.
# ---------------------------------------------------------------------------

$wsMeta.Range("A2").Value = "nomindicador"
$wsMeta.Range("B2").Value = "Porcentaje de personas de 20 años o más que no culminaron educación media superior"

$wsMeta.Range("A3").Value = "derecho"
$wsMeta.Range("B3").Value = "Educación"

$wsMeta.Range("A4").Value = "conindicador"
$wsMeta.Range("B4").Value = "No culminación de educación media superior (mayores de 20 años)"

$wsMeta.Range("A5").Value = "tipoind"
$wsMeta.Range("B5").Value = "Resultados"

$wsMeta.Range("A6").Value = "definicion"
$wsMeta.Range("B6").Value = "El indicador mide el porcentaje de personas de 20 años o más que no culminaron educación media superior."

$wsMeta.Range("A7").Value = "calculo"
$wsMeta.Range("B7").Value = "Para cada año calcular:(Cantidad de personas de 20 años o más que no culminaron la educación media superior / Cantidad de población de 20 años o más)*100"

$wsMeta.Range("A8").Value = "observaciones"
$wsMeta.Range("B8").Value = "Sin observaciones"

$wsMeta.Range("A9").Value = "cita"
$wsMeta.Range("B9").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"

$wsMeta.Range("A10").Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$wsMeta.Range("B10").Value = " "
